# Applies crypto price/volume updates per commit "Updated cryptos list" run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.299.11'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '2.569.60'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.01'
$ws.Range('E5').Value = '  +3.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.21'
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.608'
$ws.Range('E8').Value = '  +5.28%  '
$ws.Range('E9').Value = '  +4.26%  '
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('E12').Value = '  +1.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.55'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').Value = '3.029.44'
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('D15').Value = '63.253.66'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D17').Value = '2.582.55'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.37'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '343.15'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.41'
$ws.Range('E20').Value = '  +3.57%  '
$ws.Range('E21').Value = '  +1.83%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.70'
$ws.Range('E23').Value = '  +3.25%  '
$ws.Range('D24').Value = '2.688.54'
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.64'
$ws.Range('E25').Value = '  +3.57%  '
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.23'
$ws.Range('E27').Value = '  +14.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.54'
$ws.Range('E28').Value = '  +2.87%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('E31').Value = '  +8.30%  '
$ws.Range('D32').Value = '0.0₃0827'
$ws.Range('E32').Value = '  +2.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '457.64'
$ws.Range('E33').Value = '  +13.61%  '
$ws.Range('E34').Value = '  +4.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '176.81'
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.407'
$ws.Range('E36').Value = '  +3.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.28'
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.51'
$ws.Range('E38').Value = '  +4.91%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '151.76'
$ws.Range('E43').Value = '  +2.80%  '
$ws.Range('E44').Value = '  +3.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0555'
$ws.Range('E45').Value = '  +7.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.616'
$ws.Range('E46').Value = '  +3.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0984'
$ws.Range('E47').Value = '  +2.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0240'
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.49'
$ws.Range('E49').Value = '  +1.52%  '
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('E51').Value = '  -0.09%  '
